# Update "想去人数" (interested attendance count) values in column F
# across all four worksheets, per the regenerated gh-pages data output.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 1607
$ws1.Range("F10").Value = 2698
$ws1.Range("F11").Value = 2698
$ws1.Range("F13").Value = 1789
$ws1.Range("F15").Value = 290
$ws1.Range("F16").Value = 702
$ws1.Range("F17").Value = 5147
$ws1.Range("F19").Value = 84
$ws1.Range("F31").Value = 489
$ws1.Range("F34").Value = 9
$ws1.Range("F35").Value = 71
$ws1.Range("F37").Value = 58
$ws1.Range("F38").Value = 1463
$ws1.Range("F39").Value = 23
$ws1.Range("F40").Value = 1417
$ws1.Range("F41").Value = 92

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F9").Value = 117
$ws2.Range("F11").Value = 153
$ws2.Range("F18").Value = 263

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 898
$ws3.Range("F7").Value = 64

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 898
$ws4.Range("F9").Value = 64
$ws4.Range("F13").Value = 1607
$ws4.Range("F17").Value = 2698
$ws4.Range("F20").Value = 1789
$ws4.Range("F21").Value = 153
$ws4.Range("F23").Value = 290
$ws4.Range("F24").Value = 702
$ws4.Range("F25").Value = 5147
$ws4.Range("F27").Value = 84
$ws4.Range("F39").Value = 489
$ws4.Range("F42").Value = 263
$ws4.Range("F45").Value = 9
$ws4.Range("F46").Value = 71
$ws4.Range("F48").Value = 58
$ws4.Range("F49").Value = 1417
$ws4.Range("F50").Value = 92
